$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.799.25'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '''2.648.23'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''600.48'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = '''155.77'
$ws.Range('E6').Value = '  +4.09%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').Value = '''2.646.85'
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('E10').Value = '  +13.14%  '
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').Value = '''5.24'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').Value = '''0.352'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '''27.95'
$ws.Range('E14').Value = '  +2.43%  '
$ws.Range('E15').Value = '  +6.16%  '
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '''68.691.78'
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('D18').Value = '''2.645.77'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').Value = '''11.39'
$ws.Range('E19').Value = '  +3.44%  '
$ws.Range('D20').Value = '''365.78'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  +2.03%  '
$ws.Range('E24').Value = '  +4.38%  '
$ws.Range('D25').Value = '''72.92'
$ws.Range('E25').Value = '  +10.13%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '''10.08'
$ws.Range('E27').Value = '  +1.57%  '
$ws.Range('E28').Value = '  +7.86%  '
$ws.Range('D30').Value = '''583.44'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').Value = '''0.994'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '''1.42'
$ws.Range('E32').Value = '  +4.34%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''8.02'
$ws.Range('E33').Value = '  +4.99%  '
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('D35').Value = '''0.131'
$ws.Range('E35').Value = '  +5.63%  '
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +4.50%  '
$ws.Range('D38').Value = '''160.49'
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('D40').Value = '''19.36'
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('E42').Value = '  +3.62%  '
$ws.Range('D43').Value = '''2.67'
$ws.Range('E43').Value = '  +7.07%  '
$ws.Range('D44').Value = '''17.71'
$ws.Range('E44').Value = '  +5.43%  '
$ws.Range('D45').Value = '''0.0₆0322'
$ws.Range('E45').Value = '  +12.12%  '
$ws.Range('D46').Value = '''40.73'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '''156.22'
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('D49').Value = '''3.74'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').Value = '''22.07'
$ws.Range('E50').Value = '  +3.41%  '
$ws.Range('E51').Value = '  +1.61%  '
